# "40% of bridges have a length and condition, first map of the N1"
#
# Re-classify the "type" column (H) of the N1 road-network table from the
# generic "road" label into proper network-graph roles:
#   - the very first "road" row (the start of the N1)   -> "source"
#   - the very last  "road" row (the end of the N1)     -> "sink"
#   - every other row previously marked "road"          -> "link"
# Rows already marked "bridge" (or anything else) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Column H is the 8th column ("type").
$typeCol = 8

# Collect every row number whose "type" is currently "road".
$roadRows = New-Object System.Collections.ArrayList
for ($r = 2; $r -le $lastRow; $r++) {
    $val = $ws.Cells.Item($r, $typeCol).Value()
    if ($val -eq "road") {
        [void]$roadRows.Add($r)
    }
}

if ($roadRows.Count -gt 0) {
    $firstRoadRow = $roadRows[0]
    $lastRoadRow = $roadRows[$roadRows.Count - 1]

    foreach ($r in $roadRows) {
        if ($r -eq $firstRoadRow) {
            $ws.Cells.Item($r, $typeCol).Value = "source"
        } elseif ($r -eq $lastRoadRow) {
            $ws.Cells.Item($r, $typeCol).Value = "sink"
        } else {
            $ws.Cells.Item($r, $typeCol).Value = "link"
        }
    }
}
